$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, matching style of existing header row (bold, centered, bordered)
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Boolean outlier flag values for rows 2-18
$values = @(
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,1),
    @(1,1,1),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,1,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,0,0),
    @(0,1,0),
    @(0,0,0),
    @(0,0,0),
    @(1,1,1)
)

$row = 2
foreach ($rowVals in $values) {
    $ws.Cells.Item($row, 6).Value = [bool]($rowVals[0])
    $ws.Cells.Item($row, 7).Value = [bool]($rowVals[1])
    $ws.Cells.Item($row, 8).Value = [bool]($rowVals[2])
    $row++
}
